# "progress on earth :)" -- add an "initialization" sheet (colony structure
# parameters) in front of the existing sheet, which becomes "timeline" and
# gets a first timeline event row.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Sheets: insert a new sheet ("initialization") before the existing
#    one, then rename the existing sheet (by name, since Item(1) would
#    now resolve to the freshly inserted sheet) to "timeline".
# ------------------------------------------------------------------
$initSheet = $wb.Worksheets.Add()
$initSheet.Name = "initialization"

$timelineSheet = $wb.Worksheets.Item("Sheet1")
$timelineSheet.Name = "timeline"

# ------------------------------------------------------------------
# 2. "initialization" sheet: a parameter/value table for colony
#    structures, starting at row 3.
# ------------------------------------------------------------------
$headers = @("colony", "structure", "parameter", "value", "unit")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $initSheet.Cells.Item(3, $c + 1)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
}

$structures = @(
    "propelant_container",
    "booster_storage",
    "tank_storage",
    "heartofgold_storage",
    "propelant_factory",
    "booster_factory",
    "tank_factory",
    "heartofgold_factory"
)

$row = 4
foreach ($structure in $structures) {
    $initSheet.Cells.Item($row, 1).Value = "earth"
    $initSheet.Cells.Item($row, 2).Value = $structure
    $initSheet.Cells.Item($row, 3).Value = "init"
    $initSheet.Cells.Item($row, 4).Value = 0
    $initSheet.Cells.Item($row, 5).Value = "unit"
    $row++
}

# Column widths (best-fit-ish, character-width units).
$initSheet.Columns.Item(1).ColumnWidth = 5.666666666666667
$initSheet.Columns.Item(2).ColumnWidth = 16.998697916666668
$initSheet.Columns.Item(3).ColumnWidth = 8.998697916666666
$initSheet.Columns.Item(4).ColumnWidth = 4.830729166666667

$excel.ActiveWindow.Zoom = 220
$initSheet.Range("C15").Select() | Out-Null

# ------------------------------------------------------------------
# 3. "timeline" sheet: datetime/event table, existing "datetime" header
#    cell stays at A1, the rest of the header row and the first data
#    row are added.
# ------------------------------------------------------------------
$tlHeaders = @("datetime", "event", "colony", "structure", "parameter", "value", "unit")
for ($c = 0; $c -lt $tlHeaders.Length; $c++) {
    $cell = $timelineSheet.Cells.Item(1, $c + 1)
    $cell.Value = $tlHeaders[$c]
    $cell.Font.Bold = $true
}

$timelineSheet.Range("A2").NumberFormat = "m/d/yy h:mm"
$timelineSheet.Cells.Item(2, 1).Value = 45825.625
$timelineSheet.Cells.Item(2, 2).Value = "update"
$timelineSheet.Cells.Item(2, 3).Value = "earth"
$timelineSheet.Cells.Item(2, 4).Value = "propellant_factory"
$timelineSheet.Cells.Item(2, 5).Value = "rate"
$timelineSheet.Cells.Item(2, 6).Value = 10
$timelineSheet.Cells.Item(2, 7).Value = "unit/sec"

$timelineSheet.Columns.Item(1).ColumnWidth = 11.830729166666666
$timelineSheet.Columns.Item(4).ColumnWidth = 15.498697916666666

# Make "timeline" the active tab/sheet, with its own zoom + selection.
$timelineSheet.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 220
$timelineSheet.Range("G6").Select() | Out-Null
